# The workbook is re-purposed from a translation-mapping table into a
# single "Pending Review" placeholder cell: all existing data (A1:D12)
# is removed and A1 is set to the new text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe every cell currently in use on the sheet (values, not just formatting).
$ws.Cells.ClearContents()

# Write the new placeholder value into A1.
$ws.Range("A1").Value = "Pending Review"

# Make A1 the active/selected cell, matching a freshly-trimmed sheet.
$ws.Range("A1").Select() | Out-Null
